# Update output_data_dictionary.xlsx to add a new "large public green space"
# indicator (score + population percentage rows) to the Indicator estimates
# section of the data dictionary, per issue #572.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 25: "Score (/1) for access ... large public green space" ---
# This slots in right after the existing "public open space larger than 1.5
# hectares" score row, and before the GTFS public-transport score rows.
$ws.Rows("25:25").Insert()
$ws.Range("A26:D26").Copy()
$ws.Range("A25:D25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows("25:25").RowHeight = 15

$ws.Cells.Item(25, 1).Value2 = "Indicator estimates"
$ws.Cells.Item(25, 2).Value2 = "Score (/1) for access within 500m to a large public green space of at least 1 hectare in size (source: OpenStreetMap, Google Earth Engine)"
$ws.Cells.Item(25, 3).Value2 = "access_500m_large_public_green_space_score"
$ws.Cells.Item(25, 4).Value2 = "grid"

# --- Insert new row 35: "Percentage of population ... large public green space" ---
# This slots in right after the existing "public open space larger than 1.5
# hectares" population-percentage row, and before the GTFS population rows.
$ws.Rows("35:35").Insert()
$ws.Range("A34:D34").Copy()
$ws.Range("A35:D35").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows("35:35").RowHeight = 15

$ws.Cells.Item(35, 1).Value2 = "Indicator estimates"
$ws.Cells.Item(35, 2).Value2 = "Percentage of population with access within 500 m to large public green space of at least 1 hectare in size (source: OpenStreetMap, Google Earth Engine)"
$ws.Cells.Item(35, 3).Value2 = "pop_pct_access_500m_public_open_space_large_score"
$ws.Cells.Item(35, 4).Value2 = "city"

# --- Update view state to match the authored edit's cursor/selection ---
$ws.Range("A36").Select()
